# Adds the "SLA_focal" column (column L) to the harvest_day worksheet,
# mirroring the upstream commit "added sla column to harvest datashee".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell for column L
$ws.Range("L1").Value = "SLA_focal"

# New column L width (best-effort match of the authored width)
$ws.Columns.Item(12).ColumnWidth = 13

# SLA_focal values for rows 2-61
$slaValues = @(
    264.77095939999998,
    277.8731343,
    319.72162739999999,
    244.12356320000001,
    283.33035710000001,
    223.55501219999999,
    241.34086099999999,
    280.85088459999997,
    325.2873563,
    211.5696203,
    275.36694319999998,
    238.82316700000001,
    262.5930851,
    284.12444440000002,
    299.94928149999998,
    341.27966980000002,
    285.00408160000001,
    333.77483439999997,
    347.36785329999998,
    333.44978170000002,
    224.3722564,
    329.81160510000001,
    296.03155340000001,
    215.73488370000001,
    212.93056809999999,
    238.4995112,
    326.88958009999999,
    290.70403279999999,
    362.70758119999999,
    277.15582449999999,
    226.25986080000001,
    268.96243290000001,
    353.21671529999998,
    327.48695650000002,
    273.89583329999999,
    269.50784590000001,
    195.1815642,
    251.65618449999999,
    197.5708314,
    335.85629920000002,
    302.34732819999999,
    283.56313990000001,
    301.8604651,
    327.28971960000001,
    271.60000000000002,
    289.25642429999999,
    288.23715920000001,
    249.72274440000001,
    318.85462560000002,
    268.05910540000002,
    264.42801559999998,
    239.8499674,
    250.7137491,
    287.78741869999999,
    298.67562379999998,
    223.5810056,
    256.61100800000003,
    224.21785940000001,
    269.3090909,
    266.81818179999999
)

for ($i = 0; $i -lt $slaValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 12).Value = $slaValues[$i]
}

# Match the selection left behind by the edit
$ws.Range("M3").Select()
